$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update: append today's Pick 3 draw result as a new row (row 6).
# Force text formatting first so date-looking / number-looking values
# ("2025-09-22", "250922") are stored as plain text, matching the rest
# of the sheet's column data.
$ws.Range("A6:E6").NumberFormat = "@"

$ws.Range("A6").Value = "2025-09-22"
$ws.Range("B6").Value = "Pick 3"
$ws.Range("C6").Value = "250922"
$ws.Range("D6").Value = "4-9-7"
$ws.Range("E6").Value = "2025-09-22T21:37:06.648+04:00"

# Reset to the workbook's default (unstyled) cell style so the new row
# doesn't pick up a distinct style index from the temporary "@" format.
$ws.Range("A6:E6").Style = "Normal"
